# Updated cryptos list (price / 1h volume refresh, and BabyDogeCoin/EOS
# row swap) as produced by the scheduled GitHub Actions scraper run.
#
# NOTE: several "Price" (column D) values look numeric (e.g. "1.006",
# "0.5239") but must stay plain text, exactly like the original cells
# (t="inlineStr"). Setting .Value directly on such strings makes Excel
# silently coerce them into floating point numbers (losing trailing
# zeros / changing precision) and tags the cell with a "Text" style.
# To avoid both problems we temporarily force the cell's NumberFormat
# to Text ("@") before assigning the value, then reset the cell style
# back to "Normal" so the resulting cell ends up with the same (default)
# style as before, just holding a text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.519.62"
$ws.Range("E2").Value = "  +2.83%  "
$ws.Range("D3").Value = "2.123.87"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "346.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5239"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4465"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.79%  "
$ws.Range("E9").Value = "  +5.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09400"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.178"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.712"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.60%  "
$ws.Range("D14").Value = "2.131.84"
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.966"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.79%  "
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06733"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.346"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").Value = "30.558.75"
$ws.Range("E23").Value = "  +2.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.330"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("D26").Value = "2.378.25"
$ws.Range("E26").Value = "  +1.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.547"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.160"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.777"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.1062"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.880"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.293"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.44%  "
$ws.Range("E38").Value = "  +3.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06866"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.7130"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("E42").Value = "  +4.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2243"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6940"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.387"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.69%  "
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.328"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +13.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.661"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.231"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000343"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.06%  "
